# Get all V1 api tests working
#
# The "attributes" sheet describes the EMX metadata for the "ApiTestFile"
# entity. Its second attribute row was still named "fileAttr" - rename it
# to "file" so it matches the actual column name used on the "ApiTestFile"
# sheet. Also restore "attributes" as the active/selected sheet (it had
# drifted to "ApiTestFile") and refresh the remembered cell selections on
# both sheets.

$wb = $excel.ActiveWorkbook

$wsAttributes = $wb.Worksheets.Item("attributes")
$wsApiTestFile = $wb.Worksheets.Item("ApiTestFile")

# Rename the attribute from "fileAttr" to "file"
$wsAttributes.Range("A3").Value = "file"

# Update the remembered selection on the (now inactive) ApiTestFile sheet
$null = $wsApiTestFile.Range("B1").Select()

# Select a cell on "attributes" and activate it last so it becomes the
# active/selected tab of the workbook
$null = $wsAttributes.Range("C19").Select()
$wsAttributes.Activate() | Out-Null
